$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10:N10").Copy() | Out-Null
$ws.Range("A11:N11").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(11, 1).Value = 42619.894375000003
$ws.Cells.Item(11, 2).Value = -2
$ws.Cells.Item(11, 3).Value = 54
$ws.Cells.Item(11, 4).Value = 41
$ws.Cells.Item(11, 5).Value = 54
$ws.Cells.Item(11, 6).Value = 76
$ws.Cells.Item(11, 7).Value = 10593
$ws.Cells.Item(11, 8).Value = 9717
$ws.Cells.Item(11, 9).Value = 1369
$ws.Cells.Item(11, 10).Value = 259
$ws.Cells.Item(11, 11).Value = 198
$ws.Cells.Item(11, 12).Value = 5
$ws.Cells.Item(11, 13).Value = 16
$ws.Cells.Item(11, 14).Value = "Noun"
